# Weekly update: insert the newest week's price record for
# Fruta, Feria Lagunitas de Puerto Montt - Mandarina.
#
# A new data row is inserted at row 423 (pushing the previously existing
# rows 423..524 down to 424..525), and the new row is filled in with the
# latest week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 423:524 down to 424:525, leaving a blank row at 423.
$ws.Range("A423:T423").Insert()

# Fill the new row 423 with the new week's record.
$ws.Range("A423").Value = 4
$ws.Range("B423").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C423").Value = "Los Lagos"
$ws.Range("D423").Value = 45173
$ws.Range("E423").Value = 10
$ws.Range("F423").Value = "Fruta"
$ws.Range("G423").Value = 100102
$ws.Range("H423").Value = "Cítricos"
$ws.Range("I423").Value = 100102004
$ws.Range("J423").Value = "Mandarina"
$ws.Range("K423").Value = "Murcott"
$ws.Range("L423").Value = "Segunda"
$ws.Range("M423").Value = 200
$ws.Range("N423").Value = 8000
$ws.Range("O423").Value = 8000
$ws.Range("P423").Value = 8000
$ws.Range("Q423").Value = "$/bandeja 10 kilos"
$ws.Range("R423").Value = "Región de O'Higgins"
$ws.Range("S423").Value = 800
$ws.Range("T423").Value = 10
